# Update the "as of" date in the confidential disclosure banner (A7) and
# refresh the weight / percent-change figures for the model holdings table
# (D2:E3 and E4).
#
# The worksheet is protected (sheetProtection, legacy password hash), so a
# direct Range.Value assignment on a protected sheet throws. Rather than
# calling Worksheet.Unprotect (which would discard/replace the stored
# legacy password hash with a brand new one on Protect, changing parts of
# the file that are not part of this edit), we temporarily unlock just the
# cells we need to touch, make the edits, and restore their lock + number
# format state, leaving sheet-level protection completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Temporarily unlock the cells we need to edit -------------------------
$ws.Range("A7").Locked = $false
$ws.Range("D2:E3").Locked = $false
$ws.Range("E4").Locked = $false

# --- Apply the actual content changes --------------------------------------
$ws.Range("A7").Value2 = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-28 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value2 = 0.843992780349996
$ws.Range("E2").Value2 = 0.00253581843539985
$ws.Range("D3").Value2 = 0.1560072196500039
$ws.Range("E3").Value2 = 0.007499542710810214
$ws.Range("E4").Value2 = 0.003310195258715787

# --- Restore original formatting/protection on touched cells --------------
# A7 had no explicit style (default/locked); copy formats back from an
# untouched neighbor (A8) so it reverts to the same default appearance.
$ws.Range("A8").Copy()
$ws.Range("A7").PasteSpecial(-4122)

# D2:E3 and E4 used the percentage style (inherited from column D/E's
# default style) and were locked; copy formats back from the untouched D4
# cell, which carries the same original percentage style.
$ws.Range("D4").Copy()
$ws.Range("D2:E3").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("E4").PasteSpecial(-4122)
